$wb = $excel.ActiveWorkbook

# --- Update status text everywhere "Ready for handoff" appears ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# --- Narrow the "status" columns (zh-cn / de-de) on all three sheets ---
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 12.5

$wsZh.Columns.Item(3).ColumnWidth = 12.5
$wsDe.Columns.Item(3).ColumnWidth = 12.5
